# Update the "取得日時" (retrieved datetime) timestamp in column A
# for all data rows on the "ランサーズ" sheet, from 2025-10-27 18:26:06
# to 2025-10-27 18:33:35 (appended a new batch at 18:33 JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-27 18:33:35"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
